$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header / data values, entered in the same order the author typed
#     them (this drives the shared-string table order) ---
$ws.Range("F1").Value = "Cloth Category"
$ws.Range("F3").Value = "Men_Tops_Tees"
$ws.Range("F2").Value = "Women_Bottoms_Pants"
$ws.Range("G1").Value = "Size"
$ws.Range("H1").Value = "Color"
$ws.Range("H3").Value = "Red"
$ws.Range("H2").Value = "Blue"
$ws.Range("G3").Value = "M"
$ws.Range("G2").Value = 28

# --- Header row formatting: bold for the whole header row ---
$ws.Range("A1:F1").Font.Bold = $true

# --- Build the shared "text number format + left/center align" look across
#     the whole G:H block first (so G2:H3 ends up exactly on that style),
#     then bold just the header cells on top of it ---
$block = $ws.Range("G1:H3")
$block.NumberFormat = "@"
$block.HorizontalAlignment = -4131
$block.VerticalAlignment = -4108

$hdr2 = $ws.Range("G1:H1")
$hdr2.Font.Bold = $true

# --- Column widths for E, F to fit the new long header/values ---
$ws.Columns("E").ColumnWidth = 18.59
$ws.Columns("F").ColumnWidth = 20.25

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Final selection, matching the saved view state ---
$ws.Range("F6").Select()
